# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) values on the last data row
# (row 4) of the zh-cn and de-de language sheets to reflect newly
# generated report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-03-04 10:23:31"
$wsZhCn.Range("G4").Value = "2016-03-04 10:25:01"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-03-04 10:23:46"
$wsDeDe.Range("G4").Value = "2016-03-04 10:25:29"
